$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was recorded for Mango at Terminal Hortofrutícola
# Agro Chillán. It slots in chronologically right before the existing row
# for 2023-05-04 (currently row 119), so insert a fresh row there and push
# every following row down by one (the former last row, 169, becomes 170).
$ws.Rows(119).Insert()

$ws.Range("A119").Value = 7
$ws.Range("B119").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C119").Value = "Ñuble"
$ws.Range("D119").Value = 45134
$ws.Range("E119").Value = 16
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100108
$ws.Range("H119").Value = "Tropicales y subtropicales"
$ws.Range("I119").Value = 100108002
$ws.Range("J119").Value = "Mango"
$ws.Range("K119").Value = "Sin especificar"
$ws.Range("L119").Value = "Primera"
$ws.Range("M119").Value = 50
$ws.Range("N119").Value = 8000
$ws.Range("O119").Value = 8000
$ws.Range("P119").Value = 8000
$ws.Range("Q119").Value = "$/bandeja 4 kilos"
$ws.Range("R119").Value = "Brasil"
$ws.Range("S119").Value = 2000
$ws.Range("T119").Value = 4
